# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated counts (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rId1 / sheet1.xml) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 562
$wsExhibit.Range("F4").Value  = 209
$wsExhibit.Range("F7").Value  = 111
$wsExhibit.Range("F10").Value = 6891
$wsExhibit.Range("F11").Value = 243
$wsExhibit.Range("F12").Value = 383
$wsExhibit.Range("F13").Value = 3244
$wsExhibit.Range("F14").Value = 221
$wsExhibit.Range("F15").Value = 395
$wsExhibit.Range("F17").Value = 565
$wsExhibit.Range("F18").Value = 36

# --- Sheet "全部类型" (rId4 / sheet4.xml) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 562
$wsAll.Range("F6").Value  = 209
$wsAll.Range("F9").Value  = 111
$wsAll.Range("F13").Value = 6891
$wsAll.Range("F15").Value = 243
$wsAll.Range("F16").Value = 383
$wsAll.Range("F17").Value = 3244
$wsAll.Range("F18").Value = 221
$wsAll.Range("F19").Value = 395
$wsAll.Range("F21").Value = 565
$wsAll.Range("F22").Value = 36
